$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.460.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.368.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.678'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.52%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +19.10%  '

$ws.Range("E10").Value = '  +9.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.02%  '

$ws.Range("E12").Value = '  +2.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.718.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.88%  '

$ws.Range("E15").Value = '  +6.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.904'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.75%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.362.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.376.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.44%  '

$ws.Range("E19").Value = '  +6.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.93%  '

$ws.Range("E21").Value = '  +3.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '255.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.95%  '

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("E24").Value = '  -4.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.11%  '

$ws.Range("E28").Value = '  +0.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.12%  '

$ws.Range("E30").Value = '  +5.62%  '

$ws.Range("E31").Value = '  +3.23%  '

$ws.Range("E32").Value = '  +5.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0740'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.43%  '

$ws.Range("E34").Value = '  +4.72%  '

$ws.Range("E35").Value = '  +3.95%  '

$ws.Range("E36").Value = '  +8.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.43'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0272'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.07%  '

$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.89'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("E43").Value = '  +3.53%  '

$ws.Range("E44").Value = '  +3.44%  '

$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("E46").Value = '  +3.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("E48").Value = '  +12.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.442.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
